# Add a new "Greece" market tab, built as a copy of the existing
# "Croatia" tab (same layout/styles), placed right after it, and
# becomes the active/selected sheet -- matching the "Test data for
# Greece Market" commit.

$wb = $excel.ActiveWorkbook

$croatia = $wb.Worksheets.Item("Croatia")

# Mimic the "select whole sheet, then copy" UI action so the source
# sheet's selection state collapses to a full-sheet selection (as seen
# in the target Croatia sheet after the edit) instead of keeping its
# old single-cell selection.
$croatia.Range("A1:XFD1048576").Select() | Out-Null

# Copy Croatia to a new sheet placed immediately after it.
$croatia.Copy($null, $croatia)

# The newly inserted copy is now the last sheet.
$greece = $wb.Worksheets.Item($wb.Worksheets.Count)
$greece.Name = "Greece"

# Fill in the Greece-specific values (same cells that hold the
# "Croatia Market" / "NGC-xxxx" values on the source sheet).
$greece.Range("B2").Value = "Greece Market"
$greece.Range("B4").Value = "NGC-4119/T3190"

# Leave the new sheet focused/selected on its own cell, as the new tab.
$greece.Range("F17").Select() | Out-Null
